{"js": "// Replace the intro paragraph (\"You are participating ...\") with a\n// single merged run of text that also mentions the Gemini constellation\n// (was \"Perseus\"), matching the author's edit.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst newText =\n  \"You are participating in a global campaign to observe and record the \" +\n  \"faintest stars visible as a means of measuring light pollution in a \" +\n  \"given location. By locating and observing the constellation Gemini \" +\n  \"constellation in the night sky and comparing it to stellar charts, \" +\n  \"people from around the world will learn how the lights in their \" +\n  \"community contribute to light pollution. Your contributions to the \" +\n  \"online database will document the visible nighttime sky.\";\n\n// Find the target paragraph by its distinctive leading text.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text || \"\";\n  if (t.indexOf(\"You are participating\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the target paragraph.\");\n}\n\n// Clear out all of the existing (many, separately-formatted) runs, then\n// type the full replacement text back in as a single run with no\n// explicit run formatting (mirrors the author's edit, which collapsed\n// ~25 runs into one plain run).\ntarget.clear();\ntarget.insertText(newText, \"Start\");\n\nawait context.sync();\n", "ps1": "# Replace the intro paragraph (\"You are participating ...\") with a\n# single merged run of text that also mentions the Gemini constellation\n# (was \"Perseus\"), matching the author's edit.\n\n$d = $word.ActiveDocument\n\n$newText = \"You are participating in a global campaign to observe and record the faintest stars visible as a means of measuring light pollution in a given location. By locating and observing the constellation Gemini constellation in the night sky and comparing it to stellar charts, people from around the world will learn how the lights in their community contribute to light pollution. Your contributions to the online database will document the visible nighttime sky.\"\n\n# Locate the target paragraph by its distinctive leading text.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"You are participating\")\nif (-not $found) {\n    throw \"Could not locate the target paragraph.\"\n}\n\n# Expand the found hit to the whole paragraph (wdParagraph = 4), which\n# includes the trailing paragraph mark.\n[void]$rng.Expand(4)\n\n# Delete everything in the paragraph except the paragraph mark itself\n# (End - 1), which removes all of the existing (many, separately\n# formatted) runs, then type the full replacement text back in at the\n# start of the now-empty paragraph. This produces a single run with no\n# explicit run formatting, mirroring the author's edit, which collapsed\n# ~25 runs into one plain run.\n$textRange = $d.Range($rng.Start, $rng.End - 1)\n$textRange.Delete()\n\n$insertionPoint = $d.Range($rng.Start, $rng.Start)\n[void]$insertionPoint.InsertAfter($newText)\n"}
